{"js": "// The document starts with a centered date paragraph followed by a 5-column\n// practice table (\"two-digit \u00f7 one-digit\"). Every 4th row of the table holds\n// the actual problems (the rows between them are blank for student work), so\n// in document order there are exactly 26 non-empty paragraphs: the date line\n// plus the 25 \"a\u00f7b=c, d\" cells. We walk body.paragraphs in order and replace\n// the text of each non-blank paragraph with the corresponding new value \u2014\n// this is positional (not text-keyed) so it's safe even though a couple of\n// the old values (e.g. \"32\u00f75=6, 2\") repeat at different spots in the table.\nconst newValues = [\n  \"2024-11-09 Saturday\",\n  \"52\u00f75=10, 2\", \"96\u00f73=32, 0\", \"65\u00f74=16, 1\", \"75\u00f73=25, 0\", \"46\u00f72=23, 0\",\n  \"96\u00f74=24, 0\", \"13\u00f77=1, 6\", \"94\u00f76=15, 4\", \"18\u00f78=2, 2\", \"51\u00f76=8, 3\",\n  \"57\u00f73=19, 0\", \"95\u00f78=11, 7\", \"15\u00f75=3, 0\", \"41\u00f77=5, 6\", \"10\u00f72=5, 0\",\n  \"20\u00f79=2, 2\", \"86\u00f73=28, 2\", \"73\u00f72=36, 1\", \"89\u00f73=29, 2\", \"83\u00f73=27, 2\",\n  \"37\u00f79=4, 1\", \"50\u00f76=8, 2\", \"33\u00f77=4, 5\", \"57\u00f72=28, 1\", \"71\u00f78=8, 7\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\ncontext.load(paragraphs, \"text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length && idx < newValues.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text && paragraph.text.trim().length > 0) {\n    paragraph.insertText(newValues[idx], \"Replace\");\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document is a centered date line followed by a 5-column practice\n# table (\"two-digit \u00f7 one-digit\"). Every 4th table row (1, 5, 9, 13, 17 in\n# 1-based COM indexing) holds the actual problems; the rows in between are\n# blank for student work. We overwrite each of those 25 cells in place by\n# (row, column), then fix up the date line. Note: accessing $d.Tables first\n# and then indexing into $d.Paragraphs.Item(N) can resolve to the wrong\n# paragraph in this host, so the date line is updated via a Find/Replace\n# scoped to the range before the table \u2014 which is immune to that issue and\n# also avoids ambiguity, since the old cell value \"32\u00f75=6, 2\" recurs twice\n# elsewhere in the table with different replacements.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$data = @(\n  @(\"52\u00f75=10, 2\", \"96\u00f73=32, 0\", \"65\u00f74=16, 1\", \"75\u00f73=25, 0\", \"46\u00f72=23, 0\"),\n  @(\"96\u00f74=24, 0\", \"13\u00f77=1, 6\", \"94\u00f76=15, 4\", \"18\u00f78=2, 2\", \"51\u00f76=8, 3\"),\n  @(\"57\u00f73=19, 0\", \"95\u00f78=11, 7\", \"15\u00f75=3, 0\", \"41\u00f77=5, 6\", \"10\u00f72=5, 0\"),\n  @(\"20\u00f79=2, 2\", \"86\u00f73=28, 2\", \"73\u00f72=36, 1\", \"89\u00f73=29, 2\", \"83\u00f73=27, 2\"),\n  @(\"37\u00f79=4, 1\", \"50\u00f76=8, 2\", \"33\u00f77=4, 5\", \"57\u00f72=28, 1\", \"71\u00f78=8, 7\")\n)\n\n$rowIndices = @(1, 5, 9, 13, 17)\nfor ($i = 0; $i -lt $rowIndices.Length; $i++) {\n  $r = $rowIndices[$i]\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($r, $c).Range.Text = $data[$i][$c - 1]\n  }\n}\n\n$dateRange = $d.Range(0, $t.Range.Start)\n$dateRange.Find.Execute(\"2024-11-08 Friday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-11-09 Saturday\", 2)\n"}
